$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/related-observation"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Base Definition cell for the Extension row (same base URL as Metadata URL)
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/related-observation"

# Extension.value[x] Type(s) reference URL
$elements.Range("J6").Value = "Reference(http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-population-observation)
"

# Constraint(s) column for the Extension row is cleared
$elements.Range("AI2").Value = ""
